# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-12-18 10:17:31
# Applies the session-analysis refresh: summary KPIs, the "Recorded By" value
# ordering fix, newly-recorded B1D/B1E/B1F session-12 rows, and the
# recomputed Group Statistics rollups that follow from them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Top summary block (K4:L10) - Recorded/Missing session counts & percents
# ---------------------------------------------------------------------------
$ws.Range("L6").Value = 135        # Recorded Sessions
$ws.Range("L7").Value = 3          # Missing Sessions
$ws.Range("L9").Value = "42.5%"    # Coverage %
$ws.Range("L10").Value = "72.0%"   # Average Attendance %

# ---------------------------------------------------------------------------
# 2) "Recorded By" cells that listed "dnasr281@gmail.com, System" now list
#    "System, dnasr281@gmail.com" (order swapped) - every row below.
# ---------------------------------------------------------------------------
$recordedByRows = @(8,9,10,34,35,36,60,61,62,86,87,88,112,113,114,138,139,140,164,167,191,194,218,221,245,248,272,275,299,302)
foreach ($r in $recordedByRows) {
    $ws.Range("G$r").Value = "System, dnasr281@gmail.com"
}

# ---------------------------------------------------------------------------
# 3) Six sessions that moved from "Not Recorded" to "Recorded" (session 12,
#    18/12/2025, groups B1D1/B1D2/B1E1/B1E2/B1F1/B1F2). Each row's fill
#    changes from the "Not Recorded" pink to the "Recorded" green - copy the
#    format from the prior (already "Recorded") row in the same group, then
#    write the new values.
# ---------------------------------------------------------------------------
$newlyRecorded = @(
    @{ Row = 169; FormatRow = 168; RecordedBy = "dnasr281@gmail.com"; Students = "19/23" },
    @{ Row = 196; FormatRow = 195; RecordedBy = "dnasr281@gmail.com"; Students = "27/30" },
    @{ Row = 223; FormatRow = 222; RecordedBy = "dnasr281@gmail.com"; Students = "21/25" },
    @{ Row = 250; FormatRow = 249; RecordedBy = "dnasr281@gmail.com"; Students = "21/28" },
    @{ Row = 277; FormatRow = 276; RecordedBy = "dnasr281@gmail.com"; Students = "23/26" },
    @{ Row = 304; FormatRow = 303; RecordedBy = "dnasr281@gmail.com"; Students = "22/29" }
)

foreach ($item in $newlyRecorded) {
    $r = $item.Row
    $fmtR = $item.FormatRow

    $ws.Range("A${fmtR}:I${fmtR}").Copy()
    $ws.Range("A${r}:I${r}").PasteSpecial(-4122)   # xlPasteFormats

    $ws.Range("G$r").Value = $item.RecordedBy
    $ws.Range("H$r").Value = $item.Students
    $ws.Range("I$r").Value = "Recorded"
}

# ---------------------------------------------------------------------------
# 4) Group Statistics table (rows 21-26, groups B1D1/B1D2/B1E1/B1E2/B1F1/B1F2)
#    - Recorded/Missing counts and the derived Coverage %/Avg Attendance %
#    shift now that the session-12 sessions above are recorded.
# ---------------------------------------------------------------------------
$groupStats = @(
    @{ Row = 21; Recorded = 12; Missing = 0; Coverage = "44.4%"; Avg = "77.5%" },
    @{ Row = 22; Recorded = 12; Missing = 0; Coverage = "44.4%"; Avg = "75.6%" },
    @{ Row = 23; Recorded = 12; Missing = 0; Coverage = "44.4%"; Avg = "80.3%" },
    @{ Row = 24; Recorded = 11; Missing = 1; Coverage = "40.7%"; Avg = "70.1%" },
    @{ Row = 25; Recorded = 12; Missing = 0; Coverage = "44.4%"; Avg = "69.2%" },
    @{ Row = 26; Recorded = 12; Missing = 0; Coverage = "44.4%"; Avg = "60.6%" }
)

foreach ($item in $groupStats) {
    $r = $item.Row
    $ws.Range("O$r").Value = $item.Recorded
    $ws.Range("P$r").Value = $item.Missing
    $ws.Range("R$r").Value = $item.Coverage
    $ws.Range("S$r").Value = $item.Avg
}
